# fix: costo de peajes
# Adds a "Total Peajes" column (H) to the "Gastos por Unidad" sheet and
# recomputes CPK (col G) for the units whose total cost changed because of
# the newly-accounted toll expense. Also refreshes the mirrored CPK value on
# the "Top 10 Unidades Menos Eficientes" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Gastos por Unidad" ------------------------------------------
$ws = $wb.Worksheets.Item(1)

# New header H1, matching the style of the existing header row (bold,
# centered, bordered) by copying formatting from G1 then overwriting text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Total Peajes"

# Default every data row (2..404) to 0 tolls, then patch in the handful of
# units that actually carried a toll expense.
$ws.Range("H2:H404").Value = 0

$ws.Range("H2").Value = 35361
$ws.Range("H61").Value = 1367
$ws.Range("H401").Value = 16569
$ws.Range("H402").Value = 1152
$ws.Range("H403").Value = 1051

# CPK = (Gasto Combustible + Gasto Mantenimiento + Total Peajes) / Kms Totales
# Only the rows above with nonzero Kms Totales see their CPK move.
$ws.Range("G2").Value = 1091.966828865455
$ws.Range("G61").Value = 80.9619541317647

# --- Sheet 8: "Top 10 Unidades Menos Eficientes" ---------------------------
# Mirrors DC01's CPK figure; keep it in sync with the corrected value above.
$ws8 = $wb.Worksheets.Item(8)
$ws8.Range("E2").Value = 1091.966828865455
